# Update stats for 2025-08 (row 21)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = 6243
$ws.Range("D21").Value = 5625901
$ws.Range("E21").Value = 901.15345186609
$ws.Range("F21").Value = 8.366603020308983
$ws.Range("H21").Value = 28.41673035116408
